$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 16:20"

# --- Update Alemania stats (row 8) ---
$ws.Range("B8").Value = 63929
$ws.Range("C8").Value = 1494
$ws.Range("E8").Value = 54158
$ws.Range("G8").Value = 19
$ws.Range("H8").Value = 560

# --- Update Serbia stats (row 53) ---
$ws.Range("E53").Value = 727
$ws.Range("G53").Value = 3
$ws.Range("H53").Value = 16

# --- Reorder Trinidad yTobago / Mayotte / Consejo Danes para los Refugiados
#     and refresh their statistics (rows 118-120) ---
$ws.Range("A118").Value = "Trinidad yTobago"
$ws.Range("B118").Value = 82
$ws.Range("C118").Value = 4
$ws.Range("D118").Value = 1
$ws.Range("E118").Value = 78
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 3

$ws.Range("A119").Value = "Mayotte"
$ws.Range("B119").Value = 82
$ws.Range("C119").Value = 19
$ws.Range("D119").Value = 10
$ws.Range("E119").Value = 72
$ws.Range("F119").Value = 3
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 0

$ws.Range("A120").Value = "Consejo Danes para los Refugiados"
$ws.Range("B120").Value = 81
$ws.Range("C120").Value = 0
$ws.Range("D120").Value = 2
$ws.Range("E120").Value = 71
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 8
